$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("N10").Value = 64.744
$ws.Range("O10").Value = 4259
$ws.Range("N12").Value = 4299
$ws.Range("O12").Value = 62.36
